$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.053924202919006
$ws.Range("B1").Value = 1.251293182373047
$ws.Range("C1").Value = 1.642839193344116
$ws.Range("D1").Value = 3.393104553222656
$ws.Range("E1").Value = 2.305240392684937
